# Updates rows 2-21 with new sensor readings and appends new rows 22-31
# (walkingToRunning data, timestamps 2000-2900), per the "May 9th" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "walkingToRunning"
$ws.Cells.Item(2, 3).Value = -4.746825218200684
$ws.Cells.Item(2, 4).Value = -21.48864555358887
$ws.Cells.Item(2, 5).Value = -7.668013572692871
$ws.Cells.Item(2, 6).Value = 1.236815226788552
$ws.Cells.Item(2, 7).Value = 0.1414764750477494
$ws.Cells.Item(2, 8).Value = -1.8235920880804

$ws.Cells.Item(3, 1).Value = 100
$ws.Cells.Item(3, 2).Value = "walkingToRunning"
$ws.Cells.Item(3, 3).Value = 9.662276268005373
$ws.Cells.Item(3, 4).Value = -37.98822784423828
$ws.Cells.Item(3, 5).Value = -8.327095985412598
$ws.Cells.Item(3, 6).Value = 0.1625524869817767
$ws.Cells.Item(3, 7).Value = -0.4422002719431526
$ws.Cells.Item(3, 8).Value = -1.004146238825968

$ws.Cells.Item(4, 1).Value = 200
$ws.Cells.Item(4, 2).Value = "walkingToRunning"
$ws.Cells.Item(4, 3).Value = 1.045047998428345
$ws.Cells.Item(4, 4).Value = -10.61942100524902
$ws.Cells.Item(4, 5).Value = 2.765533447265625
$ws.Cells.Item(4, 6).Value = -0.740915227015286
$ws.Cells.Item(4, 7).Value = -0.7551234747400304
$ws.Cells.Item(4, 8).Value = 0.4171670004232076

$ws.Cells.Item(5, 1).Value = 300
$ws.Cells.Item(5, 2).Value = "walkingToRunning"
$ws.Cells.Item(5, 3).Value = 5.131261825561523
$ws.Cells.Item(5, 4).Value = -28.12363815307617
$ws.Cells.Item(5, 5).Value = 22.79096603393555
$ws.Cells.Item(5, 6).Value = -0.2266252960590348
$ws.Cells.Item(5, 7).Value = 0.2213748296759771
$ws.Cells.Item(5, 8).Value = 0.3229828062436459

$ws.Cells.Item(6, 1).Value = 400
$ws.Cells.Item(6, 2).Value = "walkingToRunning"
$ws.Cells.Item(6, 3).Value = -30.11330032348633
$ws.Cells.Item(6, 4).Value = -17.59586143493652
$ws.Cells.Item(6, 5).Value = -19.26617050170898
$ws.Cells.Item(6, 6).Value = 1.314678312926899
$ws.Cells.Item(6, 7).Value = 0.2296500253361486
$ws.Cells.Item(6, 8).Value = -2.454690778492302

$ws.Cells.Item(7, 1).Value = 500
$ws.Cells.Item(7, 2).Value = "walkingToRunning"
$ws.Cells.Item(7, 3).Value = 1.400394916534424
$ws.Cells.Item(7, 4).Value = -9.052268028259276
$ws.Cells.Item(7, 5).Value = -1.365690350532532
$ws.Cells.Item(7, 6).Value = 2.428599177606864
$ws.Cells.Item(7, 7).Value = -2.602694208258827
$ws.Cells.Item(7, 8).Value = -1.646401987960016

$ws.Cells.Item(8, 1).Value = 600
$ws.Cells.Item(8, 2).Value = "walkingToRunning"
$ws.Cells.Item(8, 3).Value = 29.48022842407227
$ws.Cells.Item(8, 4).Value = -27.68916893005371
$ws.Cells.Item(8, 5).Value = -8.951043128967285
$ws.Cells.Item(8, 6).Value = -6.127065071206991
$ws.Cells.Item(8, 7).Value = -1.783086014899204
$ws.Cells.Item(8, 8).Value = 3.508786747787188

$ws.Cells.Item(9, 1).Value = 700
$ws.Cells.Item(9, 2).Value = "walkingToRunning"
$ws.Cells.Item(9, 3).Value = 20.90522193908692
$ws.Cells.Item(9, 4).Value = -18.83166885375977
$ws.Cells.Item(9, 5).Value = 3.978492736816406
$ws.Cells.Item(9, 6).Value = -8.945596145478319
$ws.Cells.Item(9, 7).Value = -0.7774270279517974
$ws.Cells.Item(9, 8).Value = 0.9836971909794454

$ws.Cells.Item(10, 1).Value = 800
$ws.Cells.Item(10, 2).Value = "walkingToRunning"
$ws.Cells.Item(10, 3).Value = 6.43248987197876
$ws.Cells.Item(10, 4).Value = -11.59229469299316
$ws.Cells.Item(10, 5).Value = -0.6353058815002441
$ws.Cells.Item(10, 6).Value = -0.7768943598925588
$ws.Cells.Item(10, 7).Value = 7.095907495511296
$ws.Cells.Item(10, 8).Value = 1.043390788779353

$ws.Cells.Item(11, 1).Value = 900
$ws.Cells.Item(11, 2).Value = "walkingToRunning"
$ws.Cells.Item(11, 3).Value = 3.05394172668457
$ws.Cells.Item(11, 4).Value = -49.75492095947266
$ws.Cells.Item(11, 5).Value = 13.27557945251465
$ws.Cells.Item(11, 6).Value = 3.791030067679126
$ws.Cells.Item(11, 7).Value = 0.2994813682227342
$ws.Cells.Item(11, 8).Value = 0.2082552271494167

$ws.Cells.Item(12, 1).Value = 1000
$ws.Cells.Item(12, 2).Value = "walkingToRunning"
$ws.Cells.Item(12, 3).Value = -3.964067220687866
$ws.Cells.Item(12, 4).Value = -5.915932655334473
$ws.Cells.Item(12, 5).Value = -9.63399600982666
$ws.Cells.Item(12, 6).Value = 9.398914166633652
$ws.Cells.Item(12, 7).Value = 5.85367529281717
$ws.Cells.Item(12, 8).Value = 1.256292052073589

$ws.Cells.Item(13, 1).Value = 1100
$ws.Cells.Item(13, 2).Value = "walkingToRunning"
$ws.Cells.Item(13, 3).Value = -17.45916557312012
$ws.Cells.Item(13, 4).Value = -25.58492851257324
$ws.Cells.Item(13, 5).Value = -11.17863464355469
$ws.Cells.Item(13, 6).Value = 2.736446717322321
$ws.Cells.Item(13, 7).Value = 4.548278085443381
$ws.Cells.Item(13, 8).Value = 1.452632268533013

$ws.Cells.Item(14, 1).Value = 1200
$ws.Cells.Item(14, 2).Value = "walkingToRunning"
$ws.Cells.Item(14, 3).Value = -3.193105697631836
$ws.Cells.Item(14, 4).Value = -11.61043167114258
$ws.Cells.Item(14, 5).Value = 9.587863922119141
$ws.Cells.Item(14, 6).Value = -4.778863546864109
$ws.Cells.Item(14, 7).Value = 2.982931879182525
$ws.Cells.Item(14, 8).Value = 1.427496227997028

$ws.Cells.Item(15, 1).Value = 1300
$ws.Cells.Item(15, 2).Value = "walkingToRunning"
$ws.Cells.Item(15, 3).Value = -3.173403739929199
$ws.Cells.Item(15, 4).Value = -4.077390670776367
$ws.Cells.Item(15, 5).Value = 2.228257656097412
$ws.Cells.Item(15, 6).Value = -5.840554707887152
$ws.Cells.Item(15, 7).Value = 8.552980871390053
$ws.Cells.Item(15, 8).Value = 0.4491637055052786

$ws.Cells.Item(16, 1).Value = 1400
$ws.Cells.Item(16, 2).Value = "walkingToRunning"
$ws.Cells.Item(16, 3).Value = -78.23867034912109
$ws.Cells.Item(16, 4).Value = -53.79793548583984
$ws.Cells.Item(16, 5).Value = -33.76652908325195
$ws.Cells.Item(16, 6).Value = -1.987136993976585
$ws.Cells.Item(16, 7).Value = 6.254284772849164
$ws.Cells.Item(16, 8).Value = 0.0975735527000729

$ws.Cells.Item(17, 1).Value = 1500
$ws.Cells.Item(17, 2).Value = "walkingToRunning"
$ws.Cells.Item(17, 3).Value = 4.821199893951416
$ws.Cells.Item(17, 4).Value = -2.601359367370605
$ws.Cells.Item(17, 5).Value = -5.058528423309326
$ws.Cells.Item(17, 6).Value = 5.379759021152736
$ws.Cells.Item(17, 7).Value = 0.9356099742137876
$ws.Cells.Item(17, 8).Value = 0.3650748662601239

$ws.Cells.Item(18, 1).Value = 1600
$ws.Cells.Item(18, 2).Value = "walkingToRunning"
$ws.Cells.Item(18, 3).Value = -9.180764198303224
$ws.Cells.Item(18, 4).Value = -23.62848663330078
$ws.Cells.Item(18, 5).Value = -4.022332191467285
$ws.Cells.Item(18, 6).Value = 7.268190617395588
$ws.Cells.Item(18, 7).Value = 2.87724845140975
$ws.Cells.Item(18, 8).Value = -0.6388733892251259

$ws.Cells.Item(19, 1).Value = 1700
$ws.Cells.Item(19, 2).Value = "walkingToRunning"
$ws.Cells.Item(19, 3).Value = 1.521630764007568
$ws.Cells.Item(19, 4).Value = -8.447349548339844
$ws.Cells.Item(19, 5).Value = 13.54604339599609
$ws.Cells.Item(19, 6).Value = -1.160767073268145
$ws.Cells.Item(19, 7).Value = 1.634010728621313
$ws.Cells.Item(19, 8).Value = -0.1062023623889485

$ws.Cells.Item(20, 1).Value = 1800
$ws.Cells.Item(20, 2).Value = "walkingToRunning"
$ws.Cells.Item(20, 3).Value = 16.57039260864258
$ws.Cells.Item(20, 4).Value = -22.10472106933594
$ws.Cells.Item(20, 5).Value = 21.93498611450196
$ws.Cells.Item(20, 6).Value = -4.282147233849325
$ws.Cells.Item(20, 7).Value = 2.761744960254463
$ws.Cells.Item(20, 8).Value = 3.233038336433305

$ws.Cells.Item(21, 1).Value = 1900
$ws.Cells.Item(21, 2).Value = "walkingToRunning"
$ws.Cells.Item(21, 3).Value = -78.08035278320312
$ws.Cells.Item(21, 4).Value = -46.46374893188477
$ws.Cells.Item(21, 5).Value = -22.17394256591797
$ws.Cells.Item(21, 6).Value = -2.985469795220721
$ws.Cells.Item(21, 7).Value = 12.98016068161699
$ws.Cells.Item(21, 8).Value = -0.2386761236664393

$ws.Cells.Item(22, 1).Value = 2000
$ws.Cells.Item(22, 2).Value = "walkingToRunning"
$ws.Cells.Item(22, 3).Value = -11.79047203063965
$ws.Cells.Item(22, 4).Value = -6.283020973205566
$ws.Cells.Item(22, 5).Value = -7.190555095672607
$ws.Cells.Item(22, 6).Value = -0.1446846399875339
$ws.Cells.Item(22, 7).Value = -8.862621700526844
$ws.Cells.Item(22, 8).Value = 0.6794372058862079

$ws.Cells.Item(23, 1).Value = 2100
$ws.Cells.Item(23, 2).Value = "walkingToRunning"
$ws.Cells.Item(23, 3).Value = -2.819984912872314
$ws.Cells.Item(23, 4).Value = -18.76873397827148
$ws.Cells.Item(23, 5).Value = -7.554898738861084
$ws.Cells.Item(23, 6).Value = 10.43086243307355
$ws.Cells.Item(23, 7).Value = -2.569769633526841
$ws.Cells.Item(23, 8).Value = 0.2679488501011921

$ws.Cells.Item(24, 1).Value = 2200
$ws.Cells.Item(24, 2).Value = "walkingToRunning"
$ws.Cells.Item(24, 3).Value = 5.317728996276856
$ws.Cells.Item(24, 4).Value = -1.513343572616577
$ws.Cells.Item(24, 5).Value = -2.32539701461792
$ws.Cells.Item(24, 6).Value = -1.647490507719652
$ws.Cells.Item(24, 7).Value = 2.485880149121367
$ws.Cells.Item(24, 8).Value = 0.8528425255358673

$ws.Cells.Item(25, 1).Value = 2300
$ws.Cells.Item(25, 2).Value = "walkingToRunning"
$ws.Cells.Item(25, 3).Value = 5.113605499267578
$ws.Cells.Item(25, 4).Value = -16.18594360351562
$ws.Cells.Item(25, 5).Value = 16.00972175598145
$ws.Cells.Item(25, 6).Value = -5.678511316413173
$ws.Cells.Item(25, 7).Value = -0.6430025874384664
$ws.Cells.Item(25, 8).Value = 5.429980419329525

$ws.Cells.Item(26, 1).Value = 2400
$ws.Cells.Item(26, 2).Value = "walkingToRunning"
$ws.Cells.Item(26, 3).Value = 10.6145133972168
$ws.Cells.Item(26, 4).Value = -29.43916702270508
$ws.Cells.Item(26, 5).Value = -53.48576354980469
$ws.Cells.Item(26, 6).Value = -0.9913829153045866
$ws.Cells.Item(26, 7).Value = 14.39225333415911
$ws.Cells.Item(26, 8).Value = -1.904873143758213

$ws.Cells.Item(27, 1).Value = 2500
$ws.Cells.Item(27, 2).Value = "walkingToRunning"
$ws.Cells.Item(27, 3).Value = -6.689743518829346
$ws.Cells.Item(27, 4).Value = -5.009637832641602
$ws.Cells.Item(27, 5).Value = -2.467369556427002
$ws.Cells.Item(27, 6).Value = 1.34357805300924
$ws.Cells.Item(27, 7).Value = -3.072673494452783
$ws.Cells.Item(27, 8).Value = -1.8158495126181

$ws.Cells.Item(28, 1).Value = 2600
$ws.Cells.Item(28, 2).Value = "walkingToRunning"
$ws.Cells.Item(28, 3).Value = -15.53036594390869
$ws.Cells.Item(28, 4).Value = -35.93496704101562
$ws.Cells.Item(28, 5).Value = 9.646455764770508
$ws.Cells.Item(28, 6).Value = 4.670226090791241
$ws.Cells.Item(28, 7).Value = -2.091667162661444
$ws.Cells.Item(28, 8).Value = -4.041793617981188

$ws.Cells.Item(29, 1).Value = 2700
$ws.Cells.Item(29, 2).Value = "walkingToRunning"
$ws.Cells.Item(29, 3).Value = -8.549526214599609
$ws.Cells.Item(29, 4).Value = -4.356056213378906
$ws.Cells.Item(29, 5).Value = -4.340849876403809
$ws.Cells.Item(29, 6).Value = 0.8640256433297839
$ws.Cells.Item(29, 7).Value = 4.650747886556637
$ws.Cells.Item(29, 8).Value = -1.559826216160883

$ws.Cells.Item(30, 1).Value = 2800
$ws.Cells.Item(30, 2).Value = "walkingToRunning"
$ws.Cells.Item(30, 3).Value = 0.2473421096801757
$ws.Cells.Item(30, 4).Value = -19.03062438964844
$ws.Cells.Item(30, 5).Value = 14.19175815582275
$ws.Cells.Item(30, 6).Value = -4.644820244896485
$ws.Cells.Item(30, 7).Value = 1.258370885785854
$ws.Cells.Item(30, 8).Value = 4.646750097243193

$ws.Cells.Item(31, 1).Value = 2900
$ws.Cells.Item(31, 2).Value = "walkingToRunning"
$ws.Cells.Item(31, 3).Value = -66.04409027099609
$ws.Cells.Item(31, 4).Value = -40.48963928222656
$ws.Cells.Item(31, 5).Value = -45.51393508911133
$ws.Cells.Item(31, 6).Value = -1.00715031529102
$ws.Cells.Item(31, 7).Value = 0.9142340319046287
$ws.Cells.Item(31, 8).Value = -0.1285482684508406

